$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (column D is a text column in the source sheet)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '66.538.29'
$ws.Range("E2").Value = '  -1.34%  '
$ws.Range("D3").Value = '2.586.30'
$ws.Range("E3").Value = '  -1.96%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '584.01'
$ws.Range("E5").Value = '  -1.89%  '
$ws.Range("D6").Value = '166.56'
$ws.Range("E6").Value = '  -1.00%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -1.65%  '
$ws.Range("D9").Value = '2.585.83'
$ws.Range("E9").Value = '  -1.97%  '
$ws.Range("E10").Value = '  -4.10%  '
$ws.Range("E11").Value = '  +0.31%  '
$ws.Range("E12").Value = '  -2.00%  '
$ws.Range("E13").Value = '  -1.88%  '
$ws.Range("E14").Value = '  -4.17%  '
$ws.Range("D15").Value = '3.057.94'
$ws.Range("E15").Value = '  -1.81%  '
$ws.Range("E16").Value = '  -2.71%  '
$ws.Range("D17").Value = '66.420.28'
$ws.Range("E17").Value = '  -1.29%  '
$ws.Range("D18").Value = '2.576.71'
$ws.Range("E18").Value = '  -1.46%  '
$ws.Range("D19").Value = '11.43'
$ws.Range("E19").Value = '  -6.70%  '
$ws.Range("D20").Value = '7.72'
$ws.Range("E20").Value = '  -4.49%  '
$ws.Range("D21").Value = '352.20'
$ws.Range("E21").Value = '  -2.45%  '
$ws.Range("E22").Value = '  -3.19%  '
$ws.Range("D23").Value = '4.60'
$ws.Range("E23").Value = '  -2.39%  '
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("E25").Value = '  -4.27%  '
$ws.Range("E26").Value = '  -2.34%  '
$ws.Range("E27").Value = '  -8.66%  '
$ws.Range("D28").Value = '2.716.69'
$ws.Range("E28").Value = '  -1.90%  '
$ws.Range("D29").Value = '0.0₃0988'
$ws.Range("E29").Value = '  -3.03%  '
$ws.Range("D30").Value = '533.41'
$ws.Range("E30").Value = '  -4.28%  '
$ws.Range("D31").Value = '8.02'
$ws.Range("E32").Value = '  -3.41%  '
$ws.Range("E34").Value = '  -3.29%  '
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("E36").Value = '  -3.72%  '
$ws.Range("D37").Value = '156.90'
$ws.Range("E37").Value = '  -0.35%  '
$ws.Range("D38").Value = '18.76'
$ws.Range("E39").Value = '  -1.95%  '
$ws.Range("E40").Value = '  +1.80%  '
$ws.Range("E41").Value = '  -1.39%  '
$ws.Range("E42").Value = '  -2.14%  '
$ws.Range("E43").Value = '  +0.03%  '
$ws.Range("E44").Value = '  -2.78%  '
$ws.Range("E45").Value = '  -4.62%  '
$ws.Range("D46").Value = '149.38'
$ws.Range("E46").Value = '  -2.08%  '
$ws.Range("E47").Value = '  -3.79%  '
$ws.Range("E48").Value = '  -2.84%  '
$ws.Range("E49").Value = '  -1.86%  '
$ws.Range("D50").Value = '0.0760'
$ws.Range("E50").Value = '  -1.59%  '
$ws.Range("E51").Value = '  -1.58%  '
